$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "29.189.01", "26.01") are stored as text, matching the source data,
# then restore the default style so no stray formatting is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value2 = "29.189.01"
$ws.Range("E2").Value2 = "  +3.14%  "
$ws.Range("D3").Value2 = "1.577.40"
$ws.Range("E3").Value2 = "  +1.71%  "
$ws.Range("D4").Value2 = "0.996"
$ws.Range("E4").Value2 = "  -0.40%  "
$ws.Range("D5").Value2 = "212.03"
$ws.Range("E5").Value2 = "  +1.10%  "
$ws.Range("E6").Value2 = "  +6.70%  "
$ws.Range("D7").Value2 = "0.995"
$ws.Range("E7").Value2 = "  -0.51%  "
$ws.Range("D8").Value2 = "26.01"
$ws.Range("E8").Value2 = "  +9.97%  "
$ws.Range("E10").Value2 = "  +1.61%  "
$ws.Range("D11").Value2 = "0.0902"
$ws.Range("E11").Value2 = "  +1.36%  "
$ws.Range("D12").Value2 = "1.804.26"
$ws.Range("E12").Value2 = "  +1.78%  "
$ws.Range("D13").Value2 = "1.596.91"
$ws.Range("E13").Value2 = "  +3.06%  "
$ws.Range("D14").Value2 = "29.187.68"
$ws.Range("E14").Value2 = "  +3.15%  "
$ws.Range("B15").Value2 = "Polkadot"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value2 = "3.70"
$ws.Range("E15").Value2 = "  +2.14%  "
$ws.Range("B16").Value2 = "Polygon"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value2 = "0.522"
$ws.Range("E16").Value2 = "  +2.59%  "
$ws.Range("D17").Value2 = "62.35"
$ws.Range("E17").Value2 = "  +3.03%  "
$ws.Range("D18").Value2 = "238.62"
$ws.Range("E18").Value2 = "  +4.94%  "
$ws.Range("E19").Value2 = "  +1.79%  "
$ws.Range("D20").Value2 = "0.0₃0691"
$ws.Range("E20").Value2 = "  +2.67%  "
$ws.Range("E21").Value2 = "  -0.38%  "
$ws.Range("D22").Value2 = "3.98"
$ws.Range("E22").Value2 = "  +1.89%  "
$ws.Range("E23").Value2 = "  +4.16%  "
$ws.Range("E24").Value2 = "  +5.04%  "
$ws.Range("D25").Value2 = "153.35"
$ws.Range("E25").Value2 = "  +2.14%  "
$ws.Range("E26").Value2 = "  +4.50%  "
$ws.Range("D27").Value2 = "15.13"
$ws.Range("E27").Value2 = "  +2.59%  "
$ws.Range("E28").Value2 = "  +1.37%  "
$ws.Range("D29").Value2 = "0.996"
$ws.Range("E29").Value2 = "  -0.48%  "
$ws.Range("E30").Value2 = "  -0.17%  "
$ws.Range("E31").Value2 = "  +0.32%  "
$ws.Range("E32").Value2 = "  +1.65%  "
$ws.Range("D33").Value2 = "1.423.78"
$ws.Range("E33").Value2 = "  +2.84%  "
$ws.Range("E34").Value2 = "  +0.70%  "
$ws.Range("E36").Value2 = "  +1.67%  "
$ws.Range("D37").Value2 = "2.75"
$ws.Range("E37").Value2 = "  +6.40%  "
$ws.Range("E38").Value2 = "  -1.92%  "
$ws.Range("E39").Value2 = "  +1.01%  "
$ws.Range("D40").Value2 = "0.527"
$ws.Range("E40").Value2 = "  +3.64%  "
$ws.Range("D41").Value2 = "1.95"
$ws.Range("E41").Value2 = "  +2.10%  "
$ws.Range("D42").Value2 = "53.11"
$ws.Range("E42").Value2 = "  +26.60%  "
$ws.Range("E43").Value2 = "  -0.45%  "
$ws.Range("E44").Value2 = "  +1.56%  "
$ws.Range("E45").Value2 = "  +1.37%  "
$ws.Range("E46").Value2 = "  +4.19%  "
$ws.Range("E47").Value2 = "  -0.50%  "
$ws.Range("D48").Value2 = "1.715.39"
$ws.Range("E48").Value2 = "  +1.72%  "
$ws.Range("D49").Value2 = "0.848"
$ws.Range("E49").Value2 = "  -6.51%  "
$ws.Range("D50").Value2 = "85.75"
$ws.Range("E50").Value2 = "  +0.17%  "
$ws.Range("D51").Value2 = "0.0₆0102"
$ws.Range("E51").Value2 = "  -0.87%  "

$dRange.Style = "Normal"
